# Natmi following Dr Hou advice
# Recompute the LR-pairs table for Cadm3 (ligand) -> Cadm1 (receptor) across
# the full 3x3 sending/target cluster combination (ECs, FAPs, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "sCs")

# Per sending-cluster ligand-side stats: detected cells, rate, avg/total expr,
# avg/total derived specificity (columns E-J), shared across the 3 target rows.
$ligandStats = @{
    "ECs"  = @(2, 0.6666666666666666, 3.443291333333333, 10.329874, 0.07497468122035157, 0.07497468122035157)
    "FAPs" = @(3, 1, 39.68460733333333, 119.053822, 0.8640978924345524, 0.8640978924345523)
    "sCs"  = @(3, 1, 2.798156333333333, 8.394469000000001, 0.06092742634509613, 0.06092742634509613)
}

# Per target-cluster receptor-side stats: detected cells, rate, avg/total expr,
# avg/total derived specificity (columns K-P), shared across the 3 sending rows.
$receptorStats = @{
    "ECs"  = @(3, 1, 0.873501, 2.620503, 0.1166943280075418, 0.1166943280075418)
    "FAPs" = @(3, 1, 0.327332, 0.9819960000000001, 0.04372952953158002, 0.04372952953158002)
    "sCs"  = @(3, 1, 6.284543666666667, 18.853631, 0.8395761424608782, 0.8395761424608781)
}

# Edge-level weights/specificities (columns Q-T) for each (sending, target) pair.
$edgeStats = @{
    "ECs_ECs"   = @(3.007718422958,    27.069465806622,    0.008749120042588591, 0.008749120042588591)
    "ECs_FAPs"  = @(1.127099438722667, 10.143894948504,    0.003278607536546162, 0.003278607536546162)
    "ECs_sCs"   = @(21.63951474138822, 194.755632672494,   0.06294695364121682,  0.06294695364121682)
    "FAPs_ECs"  = @(34.664544190274,   311.9808977124659,  0.1008353228903832,   0.1008353228903832)
    "FAPs_FAPs" = @(12.99004188763467, 116.910376988712,   0.03778659430539281,  0.03778659430539281)
    "FAPs_sCs"  = @(249.3996476808535, 2244.596829127682,  0.7254759752387764,   0.7254759752387762)
    "sCs_ECs"   = @(2.444192355323,    21.997731197907,    0.007109885074569991, 0.007109885074569991)
    "sCs_FAPs"  = @(0.9159261089026668,8.243334980124002,  0.002664327689641048, 0.002664327689641048)
    "sCs_sCs"   = @(17.58513566299322, 158.266220966939,   0.05115321358088509,  0.05115321358088508)
}

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Cells.Item($row, 1).Value = $sending
        $ws.Cells.Item($row, 2).Value = "Cadm3"
        $ws.Cells.Item($row, 3).Value = "Cadm1"
        $ws.Cells.Item($row, 4).Value = $target

        $lig = $ligandStats[$sending]
        for ($i = 0; $i -lt $lig.Length; $i++) {
            $ws.Cells.Item($row, 5 + $i).Value = $lig[$i]
        }

        $rec = $receptorStats[$target]
        for ($i = 0; $i -lt $rec.Length; $i++) {
            $ws.Cells.Item($row, 11 + $i).Value = $rec[$i]
        }

        $edge = $edgeStats["$($sending)_$($target)"]
        for ($i = 0; $i -lt $edge.Length; $i++) {
            $ws.Cells.Item($row, 17 + $i).Value = $edge[$i]
        }

        $row++
    }
}
